# Update generation-mix figures on the active worksheet (rows 2-25) to
# match the revised hourly breakdown: ACTUAL_ENERGY (A), TOTAL_BCQ/SCPC/
# KSPC/KSPC1/KSPC2/EDC supply blocks (B-G) and the resulting WESM balance
# (H). Row 18's WESM figure (H18) is removed entirely (cell cleared).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 60091.7405
$ws.Range("B2").Value = 22500
$ws.Range("C2").Value = 12500
$ws.Range("D2").Value = 10000
$ws.Range("E2").Value = 5000
$ws.Range("F2").Value = 5000
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 31959.481
$ws.Range("A3").Value = 57606.156
$ws.Range("B3").Value = 22500
$ws.Range("C3").Value = 12500
$ws.Range("D3").Value = 10000
$ws.Range("E3").Value = 5000
$ws.Range("F3").Value = 5000
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 29627.31200000001
$ws.Range("A4").Value = 55826.2625
$ws.Range("B4").Value = 22500
$ws.Range("C4").Value = 12500
$ws.Range("D4").Value = 10000
$ws.Range("E4").Value = 5000
$ws.Range("F4").Value = 5000
$ws.Range("H4").Value = 27894.52499999999
$ws.Range("A5").Value = 54920.7065
$ws.Range("B5").Value = 22500
$ws.Range("C5").Value = 12500
$ws.Range("D5").Value = 10000
$ws.Range("E5").Value = 5000
$ws.Range("F5").Value = 5000
$ws.Range("H5").Value = 27069.413
$ws.Range("A6").Value = 57738.2435
$ws.Range("B6").Value = 22500
$ws.Range("D6").Value = 10000
$ws.Range("E6").Value = 5000
$ws.Range("F6").Value = 5000
$ws.Range("H6").Value = 29794.48699999999
$ws.Range("A7").Value = 62280.3545
$ws.Range("H7").Value = 34258.709
$ws.Range("A8").Value = 63002.004
$ws.Range("H8").Value = 33933.008
$ws.Range("A9").Value = 74118.9295
$ws.Range("H9").Value = 44579.859
$ws.Range("A10").Value = 89684.016
$ws.Range("B10").Value = 57000
$ws.Range("C10").Value = 25000
$ws.Range("D10").Value = 20000
$ws.Range("E10").Value = 10000
$ws.Range("F10").Value = 10000
$ws.Range("G10").Value = 12000
$ws.Range("H10").Value = 23915.03200000001
$ws.Range("A11").Value = 94141.26149999999
$ws.Range("B11").Value = 65000
$ws.Range("C11").Value = 25000
$ws.Range("G11").Value = 20000
$ws.Range("H11").Value = 14853.52299999999
$ws.Range("A12").Value = 95535.67999999999
$ws.Range("B12").Value = 65000
$ws.Range("G12").Value = 20000
$ws.Range("H12").Value = 14332.35999999999
$ws.Range("A13").Value = 96843.9975
$ws.Range("H13").Value = 15590.995
$ws.Range("A14").Value = 96558.0425
$ws.Range("H14").Value = 15364.08499999999
$ws.Range("A15").Value = 100667.975
$ws.Range("H15").Value = 19360.95000000001
$ws.Range("A16").Value = 101388.0945
$ws.Range("H16").Value = 20134.18900000001
$ws.Range("A17").Value = 80754.929
$ws.Range("H17").Value = 87.85800000000745
$ws.Range("A18").Value = 0
$ws.Range("H18").ClearContents()
$ws.Range("H19").Value = -3597.647500000006
$ws.Range("H20").Value = -454.8110000000015
$ws.Range("H21").Value = -1183.855499999998
$ws.Range("H22").Value = 138.1984999999986
$ws.Range("H23").Value = 1302.546999999991
$ws.Range("H24").Value = 2128.741999999998
$ws.Range("H25").Value = 561.2035000000033
